$wb = $excel.ActiveWorkbook

# --- Sheet 1 (currently "test 640") becomes "test 7418" ---
$sheetA = $wb.Worksheets.Item(1)
# --- Sheet 2 (currently "test 7418") becomes "test 640" ---
$sheetB = $wb.Worksheets.Item(2)

# Rename via a temporary name to avoid a name collision while swapping.
$sheetA.Name = "__swap_tmp__"
$sheetB.Name = "test 640"
$sheetA.Name = "test 7418"

# Update the data row on the sheet that is now named "test 7418"
# (the physical sheet that used to be "test 640"):
#   date 2020-11-29 -> 2020-12-05
#   odometer 30010  -> 25150
#   next service date 2021-11-29 -> 2022-12-05
# Dates are entered as literal text (not auto-converted to a date serial)
# via a text formula that is immediately converted back to a static value,
# so the cell keeps its original (unformatted) style.
$sheetA.Range("A2").Formula = '="2020-12-05"'
$sheetA.Range("A2").Copy()
$sheetA.Range("A2").PasteSpecial(-4163)

$sheetA.Range("C2").Formula = '="2022-12-05"'
$sheetA.Range("C2").Copy()
$sheetA.Range("C2").PasteSpecial(-4163)

$sheetA.Range("B2").Value = 25150.0

# Remove the data row from the sheet that is now named "test 640"
# (the physical sheet that used to be "test 7418"):
$sheetB.Rows.Item(2).Delete()
